$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.35283127929076841
$ws.Cells.Item(1, 2).Value = 0.35167003870859048
$ws.Cells.Item(2, 1).Value = -0.23628817394460633
$ws.Cells.Item(2, 2).Value = 0.23344164283753166
$ws.Cells.Item(3, 1).Value = -0.13048953651130901
$ws.Cells.Item(3, 2).Value = 0.12963781435412258
$ws.Cells.Item(4, 1).Value = -0.19163499384026927
$ws.Cells.Item(4, 2).Value = 0.19031019817672323
$ws.Cells.Item(5, 1).Value = -0.18431019844533214
$ws.Cells.Item(5, 2).Value = 0.18162246701316942
$ws.Cells.Item(6, 1).Value = -0.080626006195393174
$ws.Cells.Item(6, 2).Value = 0.080535628080339094
$ws.Cells.Item(7, 1).Value = -0.060535628415310683
$ws.Cells.Item(7, 2).Value = 0.060348283090329247
$ws.Cells.Item(8, 1).Value = -0.018260769044305292
$ws.Cells.Item(8, 2).Value = 0.018237478166689414
$ws.Cells.Item(9, 1).Value = -0.012237478451605277
$ws.Cells.Item(9, 2).Value = 0.012226525482893713
$ws.Cells.Item(10, 1).Value = -0.0062265257687172948
$ws.Cells.Item(10, 2).Value = 0.0062272638118940904
$ws.Cells.Item(11, 1).Value = -0.001727264091943681
$ws.Cells.Item(11, 2).Value = 0.0017271836650500916
$ws.Cells.Item(12, 1).Value = -0.045407982776407785
$ws.Cells.Item(12, 2).Value = 0.045156443780487443
$ws.Cells.Item(13, 1).Value = -0.039156444071348773
$ws.Cells.Item(13, 2).Value = 0.039087889773379736
$ws.Cells.Item(14, 1).Value = -0.027087890089060984
$ws.Cells.Item(14, 2).Value = 0.027054655667572369
$ws.Cells.Item(15, 1).Value = -0.021054655961103563
$ws.Cells.Item(15, 2).Value = 0.021028487596582046
$ws.Cells.Item(16, 1).Value = -0.015028487891193043
$ws.Cells.Item(16, 2).Value = 0.015004679043381763
$ws.Cells.Item(17, 1).Value = -0.0090046793394211733
$ws.Cells.Item(17, 2).Value = 0.0089999996919578606
$ws.Cells.Item(18, 1).Value = -0.10538186503261926
$ws.Cells.Item(18, 2).Value = 0.10519832243772242
$ws.Cells.Item(19, 1).Value = -0.096198322703515693
$ws.Cells.Item(19, 2).Value = 0.094728368156896448
$ws.Cells.Item(20, 1).Value = -0.085728368432334001
$ws.Cells.Item(20, 2).Value = 0.085388106488661464
$ws.Cells.Item(21, 1).Value = -0.0090044542460683274
$ws.Cells.Item(21, 2).Value = 0.0089999997219392114
$ws.Cells.Item(22, 1).Value = -0.093952106598750973
$ws.Cells.Item(22, 2).Value = 0.093637411031759044
$ws.Cells.Item(23, 1).Value = -0.084637411306905719
$ws.Cells.Item(23, 2).Value = 0.084127462019156773
$ws.Cells.Item(24, 1).Value = -0.04212746242712484
$ws.Cells.Item(24, 2).Value = 0.041999999589684478
$ws.Cells.Item(25, 1).Value = -0.094996461097295537
$ws.Cells.Item(25, 2).Value = 0.09474541318052232
$ws.Cells.Item(26, 1).Value = -0.08874541346130016
$ws.Cells.Item(26, 2).Value = 0.088422433597767025
$ws.Cells.Item(27, 1).Value = -0.082422433880346979
$ws.Cells.Item(27, 2).Value = 0.081319691255272275
$ws.Cells.Item(28, 1).Value = -0.075319691544516232
$ws.Cells.Item(28, 2).Value = 0.074555511769313298
$ws.Cells.Item(29, 1).Value = -0.062555512086541754
$ws.Cells.Item(29, 2).Value = 0.062176297307464878
$ws.Cells.Item(30, 1).Value = -0.042176297658075512
$ws.Cells.Item(30, 2).Value = 0.042020932366467534
$ws.Cells.Item(31, 1).Value = -0.027020932700953182
$ws.Cells.Item(31, 2).Value = 0.027000922735819444
$ws.Cells.Item(32, 1).Value = -0.0060009230941728831
$ws.Cells.Item(32, 2).Value = 0.0059999996997248672

# Column B width changed from 15.42578125 to 14.7109375 (stored width units).
# The ColumnWidth COM property is quantized to 1/6 character increments by
# this runtime, so 13.833333333333334 is the closest achievable setting,
# producing a stored width of 14.666666666666666 (nearest to 14.7109375).
$ws.Columns.Item(2).ColumnWidth = 13.833333333333334
